# ---------------------------------------------------------------------------
# "upload samples xlsx files"
#
# - Rename the two worksheets ("test" -> "Test1", "Sheet1" -> "Test2")
# - Duplicate the "Test1" sample data (A1:N13, plus its hyperlinks) onto the
#   previously-empty "Test2" sheet
# - Add a new trailing "new_col" column (O) of sample numbers to Test2
# - Leave Test2 as the active / selected sheet, matching the saved view
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the sheets.
$ws1.Name = "Test1"
$ws2.Name = "Test2"

# 2. Copy the whole sample table (values + shared-string reuse) from Test1
#    onto Test2 - this reproduces A1:N13 exactly, including string ids.
$ws1.Range("A1:N13").Copy()
$ws2.Range("A1").PasteSpecial()

# 3. Add the new trailing column "new_col" (O) with its sample values.
$ws2.Range("O1").Value = "new_col"
$ws2.Range("O2").Value = 12
$ws2.Range("O3").Value = 23
$ws2.Range("O4").Value = 123
$ws2.Range("O5").Value = 163.666666666667
$ws2.Range("O6").Value = 219.166666666667
$ws2.Range("O7").Value = 274.66666666666703
$ws2.Range("O8").Value = 330.16666666666703
$ws2.Range("O9").Value = 385.66666666666703
$ws2.Range("O10").Value = 441.16666666666703
$ws2.Range("O11").Value = 496.66666666666703
$ws2.Range("O12").Value = 552.16666666666697
$ws2.Range("O13").Value = 607.66666666666697

# 4. Re-create the hyperlinks that exist on Test1's C (url) and G (email)
#    columns, in the same order, against the freshly-copied Test2 cells.
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-employee-benefits-axcellent")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-business-insurance-office-insurance")
$ws2.Hyperlinks.Add($ws2.Range("G4"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("G5"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C5:C9"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-business-insurance-office-insurance", "", "", "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-business-insurance-office-insurance")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/sme-business-insurance-office-insurance")
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-business-insurance-retail-insurance")
$ws2.Hyperlinks.Add($ws2.Range("G6"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("G7"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C8"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/business-insurance-tailor-made")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/sme-business-insurance-office-insurance")
$ws2.Hyperlinks.Add($ws2.Range("C9"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/business-insurance-tailor-made")
$ws2.Hyperlinks.Add($ws2.Range("G8"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("G9"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/sme-employee-benefits-axcellent")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C10"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/sme-employee-benefits-corprotect")
$ws2.Hyperlinks.Add($ws2.Range("G10"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C11"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/sme-employee-benefits-corprotect")
$ws2.Hyperlinks.Add($ws2.Range("G11"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C12"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/en/employee-benefits-tailor-made")
$ws2.Hyperlinks.Add($ws2.Range("G12"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("G13"), "mailto:t@t.com")
$ws2.Hyperlinks.Add($ws2.Range("C13"), "https://axahk:AXA+Corp-1@hk-web-uat.paas.axa-asia.com/zh/employee-benefits-tailor-made")

# 5. Leave the view the way it was saved: Test1 showing its full table
#    selected (no single active cell highlighted), and Test2 as the
#    front-most / selected tab with P10 selected.
$ws1.Range("A1:N13").Select()
$ws2.Activate()
$ws2.Range("P10").Select()
